$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.928.85'
$ws.Range('E2').Value = '  -1.18%  '
$ws.Range('D3').Value = '2.196.09'
$ws.Range('E3').Value = '  -2.34%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '295.01'
$ws.Range('E5').Value = '  -4.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '88.58'
$ws.Range('E6').Value = '  -6.56%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.568'
$ws.Range('E7').Value = '  -0.57%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.481'
$ws.Range('E9').Value = '  -8.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.11'
$ws.Range('E10').Value = '  -7.72%  '
$ws.Range('E11').Value = '  -5.36%  '
$ws.Range('E12').Value = '  -1.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.76'
$ws.Range('E13').Value = '  -6.11%  '
$ws.Range('D14').Value = '2.529.91'
$ws.Range('E14').Value = '  -2.39%  '
$ws.Range('D15').Value = '2.261.65'
$ws.Range('E15').Value = '  -0.07%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.01'
$ws.Range('E16').Value = '  -5.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.770'
$ws.Range('E17').Value = '  -8.45%  '
$ws.Range('D18').Value = '43.546.26'
$ws.Range('E18').Value = '  -1.27%  '
$ws.Range('D19').Value = '0.0₃0887'
$ws.Range('E19').Value = '  -7.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.82'
$ws.Range('E20').Value = '  -9.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.76'
$ws.Range('E21').Value = '  -14.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '62.91'
$ws.Range('E22').Value = '  -4.61%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '230.47'
$ws.Range('E23').Value = '  -3.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.76'
$ws.Range('E24').Value = '  -12.62%  '
$ws.Range('E25').Value = '  +0.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.83'
$ws.Range('E26').Value = '  -8.85%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.16'
$ws.Range('E27').Value = '  -2.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '35.84'
$ws.Range('E28').Value = '  -6.88%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.18'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.11'
$ws.Range('E30').Value = '  -5.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '146.95'
$ws.Range('E31').Value = '  -5.00%  '
$ws.Range('E32').Value = '  -11.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.50'
$ws.Range('E33').Value = '  -5.69%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0733'
$ws.Range('E34').Value = '  -8.46%  '
$ws.Range('E35').Value = '  -3.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.86'
$ws.Range('E36').Value = '  -8.13%  '
$ws.Range('E37').Value = '  -7.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.64'
$ws.Range('E38').Value = '  -9.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0280'
$ws.Range('E39').Value = '  -8.11%  '
$ws.Range('E40').Value = '  -8.40%  '
$ws.Range('E41').Value = '  -11.61%  '
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.93'
$ws.Range('E43').Value = '  -12.56%  '
$ws.Range('D44').Value = '1.783.21'
$ws.Range('E44').Value = '  +2.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.62'
$ws.Range('E45').Value = '  +1.96%  '
$ws.Range('B46').Value = 'HuobiToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.81'
$ws.Range('E46').Value = '  +11.31%  '
$ws.Range('B47').Value = 'BitcoinSV'
$ws.Range('C47').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '72.42'
$ws.Range('E47').Value = '  -10.50%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.172'
$ws.Range('E48').Value = '  -11.19%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '13.82'
$ws.Range('E49').Value = '  +7.26%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '91.58'
$ws.Range('E50').Value = '  -8.22%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.412.95'
$ws.Range('E51').Value = '  -2.28%  '
